$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2044.6538
$ws.Range("J17").Value = 1775.0454
$ws.Range("L17").Value = 5325.1362
$ws.Range("N17").Value = -5661.1362

# Row 40
$ws.Range("H40").Value = 7642.5713
$ws.Range("J40").Value = 9999
$ws.Range("L40").Value = 9999
$ws.Range("N40").Value = -10349

# Row 74
$ws.Range("H74").Value = 15893.223
$ws.Range("J74").Value = 9999.5
$ws.Range("L74").Value = 9999.5
$ws.Range("N74").Value = -11871.5

# Row 77
$ws.Range("H77").Value = 15893.223
$ws.Range("J77").Value = 9999.5
$ws.Range("L77").Value = 49997.5
$ws.Range("N77").Value = -59357.5

# Row 113
$ws.Range("H113").Value = 4265.6665
$ws.Range("I113").Value = 4265.6665
$ws.Range("K113").Value = 4265.6665
$ws.Range("M113").Value = -1011.6665

# Row 116
$ws.Range("H116").Value = 7337.3335
$ws.Range("I116").Value = 7502.5
$ws.Range("K116").Value = 7502.5
$ws.Range("M116").Value = -4060.5

# Row 132
$ws.Range("H132").Value = 1855970
$ws.Range("I132").Value = 4205.851
$ws.Range("K132").Value = 12617.553
$ws.Range("M132").Value = -10087.553

# Row 135
$ws.Range("H135").Value = 1492.0834
$ws.Range("I135").Value = 1148.1875
$ws.Range("K135").Value = 10333.6875
$ws.Range("M135").Value = -7798.6875

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3539.4412
$ws.Range("I32").Value = 3540.6365
$ws.Range("K32").Value = 3540.6365
$ws.Range("M32").Value = -3253.6365

# Row 61
$ws.Range("H61").Value = 2261.0232
$ws.Range("I61").Value = 1681.6666
$ws.Range("K61").Value = 1681.6666
$ws.Range("M61").Value = -1469.6666

# Row 74
$ws.Range("H74").Value = 216547.84
$ws.Range("I74").Value = 329067.4
$ws.Range("J74").Value = 4010.889
$ws.Range("K74").Value = 329067.4
$ws.Range("L74").Value = 4010.889
$ws.Range("M74").Value = -328193.4
$ws.Range("N74").Value = -5758.889

# Row 77
$ws.Range("H77").Value = 216547.84
$ws.Range("I77").Value = 329067.4
$ws.Range("J77").Value = 4010.889
$ws.Range("K77").Value = 1645337
$ws.Range("L77").Value = 20054.445
$ws.Range("M77").Value = -1640969
$ws.Range("N77").Value = -28790.445

# Row 95
$ws.Range("H95").Value = 33500
$ws.Range("J95").Value = 33500
$ws.Range("L95").Value = 33500
$ws.Range("N95").Value = -38992

# Row 103
$ws.Range("H103").Value = 45000
$ws.Range("J103").Value = 45000
$ws.Range("L103").Value = 45000
$ws.Range("N103").Value = -47344

# Row 110
$ws.Range("H110").Value = 2166.4
$ws.Range("I110").Value = 2166.4
$ws.Range("K110").Value = 2166.4
$ws.Range("M110").Value = -121.4000000000001

# Row 122
$ws.Range("H122").Value = 5366.7837
$ws.Range("I122").Value = 5569.273
$ws.Range("J122").Value = 3696.25
$ws.Range("K122").Value = 16707.819
$ws.Range("L122").Value = 11088.75
$ws.Range("M122").Value = -14257.819
$ws.Range("N122").Value = -15988.75

# Row 136
$ws.Range("H136").Value = 2261.0232
$ws.Range("I136").Value = 1681.6666
$ws.Range("K136").Value = 5044.9998
$ws.Range("M136").Value = -2494.9998

$ws = $wb.Worksheets.Item("BSM")
# Row 129
$ws.Range("H129").Value = 49998
$ws.Range("J129").Value = 49998
$ws.Range("L129").Value = 49998
$ws.Range("N129").Value = -59998

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1706.3077
$ws.Range("I16").Value = 1663.3334
$ws.Range("J16").Value = 2222
$ws.Range("K16").Value = 1663.3334
$ws.Range("L16").Value = 2222
$ws.Range("M16").Value = -1376.3334
$ws.Range("N16").Value = -2796

# Row 107
$ws.Range("H107").Value = 2942035
$ws.Range("I107").Value = 5000531
$ws.Range("J107").Value = 1327
$ws.Range("K107").Value = 5000531
$ws.Range("L107").Value = 1327
$ws.Range("M107").Value = -4998611
$ws.Range("N107").Value = -5167

# Row 113
$ws.Range("H113").Value = 1706.3077
$ws.Range("I113").Value = 1663.3334
$ws.Range("J113").Value = 2222
$ws.Range("K113").Value = 1663.3334
$ws.Range("L113").Value = 2222
$ws.Range("M113").Value = 506.6666
$ws.Range("N113").Value = -6562

# Row 122
$ws.Range("H122").Value = 3847.08
$ws.Range("I122").Value = 2476.4614
$ws.Range("K122").Value = 7429.3842
$ws.Range("M122").Value = -4979.3842

# Row 132
$ws.Range("H132").Value = 16670872
$ws.Range("I132").Value = 3609.2727
$ws.Range("K132").Value = 10827.8181
$ws.Range("M132").Value = -8297.8181

$ws = $wb.Worksheets.Item("CUL")
# Row 15
$ws.Range("H15").Value = 799.2
$ws.Range("J15").Value = 1193.3334
$ws.Range("L15").Value = 3580.0002
$ws.Range("N15").Value = -3860.0002

# Row 18
$ws.Range("H18").Value = 1152.4445
$ws.Range("I18").Value = 1481.1666
$ws.Range("J18").Value = 495
$ws.Range("K18").Value = 4443.4998
$ws.Range("L18").Value = 1485
$ws.Range("M18").Value = -4274.4998
$ws.Range("N18").Value = -1823

# Row 25
$ws.Range("H25").Value = 5482.2
$ws.Range("I25").Value = 640
$ws.Range("J25").Value = 7557.4287
$ws.Range("K25").Value = 1920
$ws.Range("L25").Value = 22672.2861
$ws.Range("M25").Value = -1751
$ws.Range("N25").Value = -23010.2861

# Row 30
$ws.Range("H30").Value = 5482.2
$ws.Range("I30").Value = 640
$ws.Range("J30").Value = 7557.4287
$ws.Range("K30").Value = 1920
$ws.Range("L30").Value = 22672.2861
$ws.Range("M30").Value = -1818
$ws.Range("N30").Value = -22876.2861

# Row 111
$ws.Range("H111").Value = 167275.67
$ws.Range("I111").Value = 167275.67
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 501827.01
$ws.Range("L111").Value = 0
$ws.Range("N111").Value = -498760.01

# Row 112
$ws.Range("H112").Value = 2508.8572
$ws.Range("I112").Value = 1760.3334
$ws.Range("K112").Value = 5281.0002
$ws.Range("M112").Value = -4173.0002

# Row 126
$ws.Range("H126").Value = 6515
$ws.Range("J126").Value = 7000
$ws.Range("L126").Value = 21000
$ws.Range("N126").Value = -30880

# Row 136
$ws.Range("H136").Value = 1960.3334
$ws.Range("I136").Value = 948.875
$ws.Range("K136").Value = 2846.625
$ws.Range("M136").Value = 2253.375

# Row 139
$ws.Range("H139").Value = 5544.1763
$ws.Range("I139").Value = 12270.3
$ws.Range("J139").Value = 2741.625
$ws.Range("K139").Value = 36810.89999999999
$ws.Range("L139").Value = 8224.875
$ws.Range("M139").Value = -31670.89999999999
$ws.Range("N139").Value = -18504.875

$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 55107
$ws.Range("I57").Value = 49658
$ws.Range("K57").Value = 49658
$ws.Range("M57").Value = -48838

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 787.25
$ws.Range("I16").Value = 661.6667
$ws.Range("J16").Value = 1164
$ws.Range("K16").Value = 661.6667
$ws.Range("L16").Value = 1164
$ws.Range("M16").Value = -491.6667
$ws.Range("N16").Value = -1504

# Row 46
$ws.Range("H46").Value = 2642.5715
$ws.Range("J46").Value = 2500
$ws.Range("L46").Value = 2500
$ws.Range("N46").Value = -2876

# Row 138
$ws.Range("H138").Value = 79938
$ws.Range("J138").Value = 79938
$ws.Range("L138").Value = 79938
$ws.Range("N138").Value = -90218

$ws = $wb.Worksheets.Item("WVR")
# Row 92
$ws.Range("H92").Value = 150000
$ws.Range("J92").Value = 150000
$ws.Range("L92").Value = 150000
$ws.Range("N92").Value = -154992

# Row 98
$ws.Range("H98").Value = 45000
$ws.Range("J98").Value = 45000
$ws.Range("L98").Value = 45000
$ws.Range("N98").Value = -50990

# Row 136
$ws.Range("H136").Value = 23857.967
$ws.Range("I136").Value = 24355.26
$ws.Range("K136").Value = 73065.78
$ws.Range("M136").Value = -70515.78
